{"js": "// Update the date paragraph and the 25 division-problem answers in the\n// practice table. Every edit only swaps the text inside an existing\n// <w:t> run (formatting, paragraph/table structure stay untouched), so we\n// locate each run's current text with Search and replace just that range.\n\nasync function replaceFirst(rangeLikeObject, oldText, newText) {\n  const results = rangeLikeObject.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\n// 1) Title / date line.\nawait replaceFirst(context.document.body, \"2024-11-21 Thursday\", \"2024-11-22 Friday\");\nawait context.sync();\n\n// 2) Table of division problems. The table has 20 rows; data lives in\n// rows 0, 4, 8, 12, 16 (5 columns each), the rows in between are blank\n// spacer rows.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst dataRows = [0, 4, 8, 12, 16];\nconst values = [\n  [\"63\u00f78=7, 7\", \"98\u00f74=24, 2\"],\n  [\"90\u00f79=10, 0\", \"33\u00f73=11, 0\"],\n  [\"97\u00f72=48, 1\", \"15\u00f79=1, 6\"],\n  [\"64\u00f73=21, 1\", \"16\u00f72=8, 0\"],\n  [\"23\u00f73=7, 2\", \"80\u00f72=40, 0\"],\n  [\"45\u00f75=9, 0\", \"80\u00f73=26, 2\"],\n  [\"21\u00f79=2, 3\", \"82\u00f79=9, 1\"],\n  [\"41\u00f72=20, 1\", \"93\u00f75=18, 3\"],\n  [\"32\u00f75=6, 2\", \"77\u00f75=15, 2\"],\n  [\"49\u00f79=5, 4\", \"67\u00f78=8, 3\"],\n  [\"83\u00f77=11, 6\", \"27\u00f73=9, 0\"],\n  [\"15\u00f72=7, 1\", \"78\u00f77=11, 1\"],\n  [\"63\u00f78=7, 7\", \"29\u00f73=9, 2\"],\n  [\"14\u00f72=7, 0\", \"67\u00f78=8, 3\"],\n  [\"66\u00f76=11, 0\", \"15\u00f78=1, 7\"],\n  [\"39\u00f72=19, 1\", \"96\u00f75=19, 1\"],\n  [\"47\u00f78=5, 7\", \"35\u00f74=8, 3\"],\n  [\"69\u00f77=9, 6\", \"19\u00f77=2, 5\"],\n  [\"48\u00f75=9, 3\", \"45\u00f79=5, 0\"],\n  [\"66\u00f73=22, 0\", \"77\u00f73=25, 2\"],\n  [\"12\u00f72=6, 0\", \"42\u00f72=21, 0\"],\n  [\"69\u00f74=17, 1\", \"37\u00f72=18, 1\"],\n  [\"63\u00f77=9, 0\", \"28\u00f73=9, 1\"],\n  [\"53\u00f72=26, 1\", \"90\u00f72=45, 0\"],\n  [\"43\u00f75=8, 3\", \"60\u00f79=6, 6\"],\n];\n\nlet i = 0;\nfor (const row of dataRows) {\n  for (let col = 0; col < 5; col++) {\n    const [oldText, newText] = values[i++];\n    const cell = table.getCell(row, col);\n    await replaceFirst(cell.body, oldText, newText);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date paragraph and the 25 division-problem answers in the\n# practice table. Each edit only swaps the text inside an existing run\n# (formatting / paragraph / table structure are left untouched), so every\n# replacement is scoped to a single Range (document title, or one table\n# cell) and uses wdReplaceOne + wdFindStop so a value that happens to\n# repeat elsewhere in the table is not touched.\n\n$d = $word.ActiveDocument\n\nfunction Replace-InRange($range, $oldText, $newText) {\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    # MatchCase=$true, MatchWholeWord=$false, MatchWildcards=$false,\n    # MatchSoundsLike=$false, MatchAllWordForms=$false, Forward=$true,\n    # Wrap=wdFindStop(0), Format=$false, Replace=wdReplaceOne(1)\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 0, $false, $newText, 1) | Out-Null\n}\n\n# 1) Title / date line.\nReplace-InRange $d.Paragraphs(1).Range \"2024-11-21 Thursday\" \"2024-11-22 Friday\"\n\n# 2) Table of division problems. The table has 20 rows; data lives in\n# rows 1, 5, 9, 13, 17 (5 columns each, 1-based), the rows in between are\n# blank spacer rows.\n$table = $d.Tables.Item(1)\n\n$dataRows = @(1, 5, 9, 13, 17)\n$values = @(\n    @(\"63\u00f78=7, 7\", \"98\u00f74=24, 2\"),\n    @(\"90\u00f79=10, 0\", \"33\u00f73=11, 0\"),\n    @(\"97\u00f72=48, 1\", \"15\u00f79=1, 6\"),\n    @(\"64\u00f73=21, 1\", \"16\u00f72=8, 0\"),\n    @(\"23\u00f73=7, 2\", \"80\u00f72=40, 0\"),\n    @(\"45\u00f75=9, 0\", \"80\u00f73=26, 2\"),\n    @(\"21\u00f79=2, 3\", \"82\u00f79=9, 1\"),\n    @(\"41\u00f72=20, 1\", \"93\u00f75=18, 3\"),\n    @(\"32\u00f75=6, 2\", \"77\u00f75=15, 2\"),\n    @(\"49\u00f79=5, 4\", \"67\u00f78=8, 3\"),\n    @(\"83\u00f77=11, 6\", \"27\u00f73=9, 0\"),\n    @(\"15\u00f72=7, 1\", \"78\u00f77=11, 1\"),\n    @(\"63\u00f78=7, 7\", \"29\u00f73=9, 2\"),\n    @(\"14\u00f72=7, 0\", \"67\u00f78=8, 3\"),\n    @(\"66\u00f76=11, 0\", \"15\u00f78=1, 7\"),\n    @(\"39\u00f72=19, 1\", \"96\u00f75=19, 1\"),\n    @(\"47\u00f78=5, 7\", \"35\u00f74=8, 3\"),\n    @(\"69\u00f77=9, 6\", \"19\u00f77=2, 5\"),\n    @(\"48\u00f75=9, 3\", \"45\u00f79=5, 0\"),\n    @(\"66\u00f73=22, 0\", \"77\u00f73=25, 2\"),\n    @(\"12\u00f72=6, 0\", \"42\u00f72=21, 0\"),\n    @(\"69\u00f74=17, 1\", \"37\u00f72=18, 1\"),\n    @(\"63\u00f77=9, 0\", \"28\u00f73=9, 1\"),\n    @(\"53\u00f72=26, 1\", \"90\u00f72=45, 0\"),\n    @(\"43\u00f75=8, 3\", \"60\u00f79=6, 6\")\n)\n\n$i = 0\nforeach ($row in $dataRows) {\n    for ($col = 1; $col -le 5; $col++) {\n        $pair = $values[$i]\n        $i++\n        $cell = $table.Cell($row, $col)\n        Replace-InRange $cell.Range $pair[0] $pair[1]\n    }\n}\n"}
